# Update simulated performance scores (perform_score) and rankings (ranking)
# to reflect the new simulated energy/nutrient recoveries and user cost.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: perform_score ---
$ws1 = $wb.Worksheets.Item("perform_score")

$performScore = @{
    2  = @(0.889137668396077, 0.3268096415522836, 0.6731903584477165)
    3  = @(0.1739569652645787, 0.42128665932963, 0.7156052394754749)
    4  = @(0.6702681412954919, 0.8033328484168533, 0.2335580059493088)
    5  = @(0.9523697074025108, 1, 0)
    6  = @(0, 0.6073985986653843, 0.6279606252811276)
    7  = @(0.682597857405016, 0.7972474240278519, 0.2422737756273134)
    8  = @(0.6875180883654941, 0.7637253843241474, 0.2687251509939895)
    9  = @(0.6186359537253805, 0.7255526997905448, 0.3100163819947988)
    10 = @(0.6760100084158057, 0.8007675067080211, 0.2371207849508889)
    11 = @(0.8361985665374894, 0.8995309989365531, 0.1203088025558765)
    12 = @(0.4889244819024822, 0.7040326301430895, 0.387963706364472)
}

foreach ($row in $performScore.Keys) {
    $vals = $performScore[$row]
    $ws1.Cells.Item($row, 3).Value = $vals[0]
    $ws1.Cells.Item($row, 4).Value = $vals[1]
    $ws1.Cells.Item($row, 5).Value = $vals[2]
}

# --- Sheet 2: ranking ---
$ws2 = $wb.Worksheets.Item("ranking")

$ranking = @{
    2  = @(3, 1, 2)
    3  = @(1, 2, 3)
    4  = @(2, 3, 1)
    5  = @(2, 3, 1)
    6  = @(1, 2, 3)
    7  = @(2, 3, 1)
    8  = @(2, 3, 1)
    9  = @(2, 3, 1)
    10 = @(2, 3, 1)
    11 = @(2, 3, 1)
    12 = @(2, 3, 1)
}

foreach ($row in $ranking.Keys) {
    $vals = $ranking[$row]
    $ws2.Cells.Item($row, 3).Value = $vals[0]
    $ws2.Cells.Item($row, 4).Value = $vals[1]
    $ws2.Cells.Item($row, 5).Value = $vals[2]
}
